$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 84, shifting existing rows 84-193 down to 85-194
$ws.Rows.Item(84).EntireRow.Insert()

# Fill in the new row 84 with the weekly data entry
$ws.Range("A84").Value = 8
$ws.Range("B84").Value = "Terminal La Palmera de La Serena"
$ws.Range("C84").Value = "Coquimbo"
$ws.Range("D84").Value = 44482
$ws.Range("E84").Value = 4
$ws.Range("F84").Value = 100112032
$ws.Range("G84").Value = "Zapallo italiano"
$ws.Range("H84").Value = "Sin especificar"
$ws.Range("I84").Value = "Primera"
$ws.Range("J84").Value = 500
$ws.Range("K84").Value = 13000
$ws.Range("L84").Value = 14000
$ws.Range("M84").Value = 13500
$ws.Range("N84").Value = "$/caja 70 unidades"
$ws.Range("O84").Value = "Provincia de Limarí"
$ws.Range("P84").Value = 193
$ws.Range("Q84").Value = 70
$ws.Range("R84").Value = "Hortaliza"
